$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Inventory")

# Add new row 4 labels "System Provided" over columns K and L
$ws.Range("K4").Value = "System Provided"
$ws.Range("L4").Value = "System Provided"

# Apply same style (left/top aligned) as the header label cells in row 5 (A5:J5)
$ws.Range("K4:L4").HorizontalAlignment = -4131
$ws.Range("K4:L4").VerticalAlignment = -4160

# Adjust column widths for K and L to match new merged sizing
$ws.Range("K1").EntireColumn.ColumnWidth = 16
$ws.Range("L1").EntireColumn.ColumnWidth = 16

# Update view: scroll so column D is the top-left visible column, and select L4
$ws.Application.ActiveWindow.ScrollColumn = 4
$ws.Range("L4").Select()

# Set page orientation to portrait
$ws.PageSetup.Orientation = 1
